# Minor typo fixes in agent instructions (Agent Instructions sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agent Instructions")

# --- CS633_Agent row (row 2): drop the redundant "You are the CS633 Agent, " lead-in ---
$prefixCs633 = "You are the CS633 Agent, "
$d2 = $ws.Range("D2").Value2
if ($d2.StartsWith($prefixCs633)) {
    $ws.Range("D2").Value2 = "You are " + $d2.Substring($prefixCs633.Length)
}

# --- Career_Agent row (row 3): drop the redundant "You are the Career Agent, " lead-in ---
$prefixCareer = "You are the Career Agent, "
$d3 = $ws.Range("D3").Value2
if ($d3.StartsWith($prefixCareer)) {
    $ws.Range("D3").Value2 = "You are " + $d3.Substring($prefixCareer.Length)
}

# --- Scheduling_Agent row (row 5): drop the redundant "You are the Scheduling Agent, " lead-in ---
$prefixSched = "You are the Scheduling Agent, "
$d5 = $ws.Range("D5").Value2
if ($d5.StartsWith($prefixSched)) {
    $ws.Range("D5").Value2 = "You are " + $d5.Substring($prefixSched.Length)
}

# --- Advisor_Agent row (row 6): broaden the declined-topics sentence ---
$oldPhrase = "Computer Science department of Boston"
$newPhrase = "Computer Science, CIS or CS at Boston"
$d6 = $ws.Range("D6").Value2
$ws.Range("D6").Value2 = $d6.Replace($oldPhrase, $newPhrase)

# --- Row heights shrink now that rows 3 & 5 hold shorter text (auto-wrap rows) ---
$ws.Rows(3).RowHeight = 304
$ws.Rows(5).RowHeight = 192

# --- Restore the editor's last selection ---
$ws.Range("D7").Select()
